# Updates cryptos list data cells (Price / Volume(1h) columns, plus two
# row-47/48 coin-name+link swaps) to match the latest scrape.
#
# Cells are written as text (NumberFormat "@") so values such as
# "56.105.28" or "10.11" are preserved verbatim instead of Excel
# auto-coercing them into numbers; Style is then reset back to "Normal"
# so no residual formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '56.105.28'
Set-TextValue 'E2' '  +6.69%  '
Set-TextValue 'D3' '2.473.25'
Set-TextValue 'E3' '  +4.99%  '
Set-TextValue 'E4' '  +0.00%  '
Set-TextValue 'D5' '488.32'
Set-TextValue 'E5' '  +7.85%  '
Set-TextValue 'D6' '145.77'
Set-TextValue 'E6' '  +13.99%  '
Set-TextValue 'E7' '  -0.03%  '
Set-TextValue 'E8' '  +8.15%  '
Set-TextValue 'D9' '2.486.92'
Set-TextValue 'E9' '  +4.81%  '
Set-TextValue 'E10' '  +11.15%  '
Set-TextValue 'D11' '0.0977'
Set-TextValue 'E11' '  +4.92%  '
Set-TextValue 'E12' '  +7.68%  '
Set-TextValue 'E13' '  +2.17%  '
Set-TextValue 'D14' '2.905.46'
Set-TextValue 'E14' '  +4.93%  '
Set-TextValue 'D15' '56.226.71'
Set-TextValue 'E15' '  +6.55%  '
Set-TextValue 'E16' '  +9.69%  '
Set-TextValue 'E17' '  +6.71%  '
Set-TextValue 'D18' '2.485.96'
Set-TextValue 'E18' '  +4.65%  '
Set-TextValue 'D19' '4.55'
Set-TextValue 'E19' '  +11.26%  '
Set-TextValue 'D20' '10.11'
Set-TextValue 'E20' '  +10.24%  '
Set-TextValue 'D21' '317.83'
Set-TextValue 'E21' '  +5.20%  '
Set-TextValue 'D22' '0.997'
Set-TextValue 'E22' '  +0.44%  '
Set-TextValue 'D23' '5.80'
Set-TextValue 'E23' '  +11.01%  '
Set-TextValue 'D24' '58.42'
Set-TextValue 'E24' '  +6.21%  '
Set-TextValue 'E25' '  +9.01%  '
Set-TextValue 'D26' '0.163'
Set-TextValue 'E26' '  +10.27%  '
Set-TextValue 'D28' '2.589.96'
Set-TextValue 'E28' '  +4.59%  '
Set-TextValue 'E29' '  +9.60%  '
Set-TextValue 'D30' '0.0₃0785'
Set-TextValue 'E30' '  +11.23%  '
Set-TextValue 'D31' '1.00'
Set-TextValue 'E31' '  +0.20%  '
Set-TextValue 'D32' '148.87'
Set-TextValue 'E32' '  +3.44%  '
Set-TextValue 'D33' '18.17'
Set-TextValue 'E33' '  +4.82%  '
Set-TextValue 'E34' '  +7.90%  '
Set-TextValue 'E35' '  +6.38%  '
Set-TextValue 'E36' '  +10.50%  '
Set-TextValue 'D37' '3.71'
Set-TextValue 'E37' '  +7.95%  '
Set-TextValue 'D38' '0.856'
Set-TextValue 'E38' '  +10.08%  '
Set-TextValue 'D39' '34.07'
Set-TextValue 'E39' '  +4.21%  '
Set-TextValue 'D40' '3.50'
Set-TextValue 'E40' '  +9.36%  '
Set-TextValue 'D41' '0.0560'
Set-TextValue 'E41' '  +8.85%  '
Set-TextValue 'D42' '0.993'
Set-TextValue 'E42' '  -0.34%  '
Set-TextValue 'D43' '0.608'
Set-TextValue 'E43' '  +4.28%  '
Set-TextValue 'E44' '  +10.34%  '
Set-TextValue 'D45' '4.78'
Set-TextValue 'E45' '  +16.90%  '
Set-TextValue 'D46' '0.0917'
Set-TextValue 'E46' '  +7.80%  '
Set-TextValue 'B47' 'VeChain'
Set-TextValue 'C47' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D47' '0.0229'
Set-TextValue 'E47' '  +8.26%  '
Set-TextValue 'B48' 'WhiteBITCoin'
Set-TextValue 'C48' 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 'D48' '10.18'
Set-TextValue 'E48' '  +0.53%  '
Set-TextValue 'D49' '257.08'
Set-TextValue 'E49' '  +20.67%  '
Set-TextValue 'D50' '1.893.78'
Set-TextValue 'E50' '  -1.07%  '
Set-TextValue 'D51' '17.59'
Set-TextValue 'E51' '  +9.13%  '
